$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J (I0, IF), matching the existing
# bold / centered / bordered header style used by the other headers (s="1").
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$hdr = $ws.Range("I1:J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous
$hdr.Borders.Weight = 2            # xlThin

# New numeric data for columns I (I0) and J (IF), rows 2-20.
$data = @(
    @(7,7),
    @(7,8),
    @(8,8),
    @(5,7),
    @(8,8),
    @(6,8),
    @(6,6),
    @(5,5),
    @(6,7),
    @(9,9),
    @(7,7),
    @(9,9),
    @(9,9),
    @(5,5),
    @(6,6),
    @(8,9),
    @(8,8),
    @(8,8),
    @(9,9)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}

Write-Output "done"
